$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "35.531.87"
$ws.Cells.Item(2, 5).Value = "  -0.19%  "

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "1.890.21"
$ws.Cells.Item(3, 5).Value = "  -0.31%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.02%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "246.05"
$ws.Cells.Item(5, 5).Value = "  -0.68%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.690"
$ws.Cells.Item(6, 5).Value = "  -0.47%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  +0.03%  "

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "43.08"
$ws.Cells.Item(8, 5).Value = "  -1.83%  "

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "56.57"
$ws.Cells.Item(9, 5).Value = "  +8.78%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.356"
$ws.Cells.Item(10, 5).Value = "  +0.83%  "

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0752"
$ws.Cells.Item(11, 5).Value = "  +1.34%  "

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.0983"
$ws.Cells.Item(12, 5).Value = "  +1.31%  "

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "14.57"
$ws.Cells.Item(13, 5).Value = "  +10.97%  "

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "0.792"
$ws.Cells.Item(14, 5).Value = "  +8.26%  "

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "2.170.54"
$ws.Cells.Item(15, 5).Value = "  -0.06%  "

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "5.02"
$ws.Cells.Item(16, 5).Value = "  +1.25%  "

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "1.895.48"
$ws.Cells.Item(17, 5).Value = "  -1.79%  "

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "35.535.11"
$ws.Cells.Item(18, 5).Value = "  -0.05%  "

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "73.43"
$ws.Cells.Item(19, 5).Value = "  -0.55%  "

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "0.0₃0829"
$ws.Cells.Item(20, 5).Value = "  +0.35%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "246.37"
$ws.Cells.Item(21, 5).Value = "  -0.38%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "13.02"
$ws.Cells.Item(22, 5).Value = "  +0.92%  "

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "5.17"
$ws.Cells.Item(23, 5).Value = "  +4.05%  "

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "2.66"
$ws.Cells.Item(24, 5).Value = "  +3.87%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  +0.02%  "

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "2.14"
$ws.Cells.Item(26, 5).Value = "  -2.42%  "

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "166.15"
$ws.Cells.Item(27, 5).Value = "  +0.11%  "

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "8.64"
$ws.Cells.Item(28, 5).Value = "  +1.27%  "

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "18.34"
$ws.Cells.Item(29, 5).Value = "  -0.54%  "

# Row 30
$ws.Cells.Item(30, 5).Value = "  -0.14%  "

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "4.41"
$ws.Cells.Item(31, 5).Value = "  +3.66%  "

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.0608"
$ws.Cells.Item(32, 5).Value = "  +4.32%  "

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "4.26"
$ws.Cells.Item(33, 5).Value = "  +0.59%  "

# Row 34
$ws.Cells.Item(34, 5).Value = "  +19.55%  "

# Row 35
$ws.Cells.Item(35, 5).Value = "  +0.04%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  -16.36%  "

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.852"
$ws.Cells.Item(37, 5).Value = "  -0.18%  "

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.0748"
$ws.Cells.Item(38, 5).Value = "  +9.37%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "1.94"
$ws.Cells.Item(39, 5).Value = "  -3.69%  "

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.0229"
$ws.Cells.Item(40, 5).Value = "  +7.17%  "

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "99.04"
$ws.Cells.Item(41, 5).Value = "  +1.07%  "

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "16.97"
$ws.Cells.Item(42, 5).Value = "  -1.32%  "

# Row 43
$ws.Cells.Item(43, 2).Value = "ARBITRUM"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "1.09"
$ws.Cells.Item(43, 5).Value = "  -0.63%  "

# Row 44
$ws.Cells.Item(44, 2).Value = "Gas"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "14.38"
$ws.Cells.Item(44, 5).Value = "  +17.99%  "

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "1.313.57"
$ws.Cells.Item(45, 5).Value = "  +1.01%  "

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "2.34"
$ws.Cells.Item(46, 5).Value = "  -1.37%  "

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.0812"
$ws.Cells.Item(47, 5).Value = "  -0.11%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  -0.07%  "

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "2.73"
$ws.Cells.Item(49, 5).Value = "  -0.60%  "

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "6.37"
$ws.Cells.Item(50, 5).Value = "  +0.11%  "

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "42.49"
$ws.Cells.Item(51, 5).Value = "  -2.05%  "
